$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "patientId" row into the resourcesInfo table (2nd table in
#    the document), right before the existing "centerName" row.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(2)

# Locate the "centerName" row so the new row can be added immediately before it.
$beforeRow = $null
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $cellText = $t.Cell($i, 1).Range.Text
    if ($cellText -like "centerName*") {
        $beforeRow = $t.Rows.Item($i)
        break
    }
}

$newRow = $t.Rows.Add($beforeRow)

$newRow.Cells.Item(1).Range.Text = "patientId"
$newRow.Cells.Item(2).Range.Text = "ID partagé du patient transporté"
$newRow.Cells.Item(3).Range.Text = "string"
$newRow.Cells.Item(4).Range.Text = "0..1"
$newRow.Cells.Item(5).Range.Text = "Identifiant partagé du patient qui est transporté. Ce n'est à remplir que lorsque l'on sait quel vecteur transporte quel patient. " + [char]11 + "Il est valorisé comme suit lors de sa création : " + [char]11 + "{OrgId émetteur}.patient.{n°patient unique dans le système émetteur}" + [char]11 + [char]11 + "OU, si un n°patient unique n'existe pas dans le système émetteur :" + [char]11 + "{ID émetteur}.{senderCaseId}.patient.{numéro d’ordre chronologique au dossier}"
$newRow.Cells.Item(6).Range.Text = "fr.health.samu440.patient.P23AZ59"

# ---------------------------------------------------------------------------
# 2) Rename the nomenclature references to the new HubSante.* naming scheme.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("(NOMENCLATURE: SI-SAMU-TYPE_VECTEUR)", $true, $false, $false, $false, $false, $true, 1, $false, "(NOMENCLATURE: HubSante.typeVecteur)", 2) | Out-Null

$d.Content.Find.Execute("(NOMENCLATURE: SI-SAMU-NIVSOIN)", $true, $false, $false, $false, $false, $true, 1, $false, "(NOMENCLATURE: HubSante.typePEC)", 2) | Out-Null

$d.Content.Find.Execute("(NOMENCLATURE: SI-SAMU-STATUS_VECTEUR)", $true, $false, $false, $false, $false, $true, 1, $false, "(NOMENCLATURE: HubSante.statutVecteur)", 2) | Out-Null

$d.Content.Find.Execute("(NOMENCLATURE: ENUM-CONTACT_Type)", $true, $false, $false, $false, $false, $true, 1, $false, "(NOMENCLATURE: HubSante.typeCom)", 2) | Out-Null
